$d = $word.ActiveDocument

function Insert-XmlOverRange($rng, $bodyXml) {
    $pkg = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p>' + $bodyXml + '</w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
    $rng.InsertXML($pkg)
}

# Change 1: In "OnInitDialog" -> wrap OnInitDialog with spellStart/spellEnd
$r1 = $d.Content
$r1.Find.Execute("In “OnInitDialog”", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$rng1 = $d.Range($r1.Start, $r1.End)
$body1 = '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>In “</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>OnInitDialog</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>”</w:t></w:r>'
Insert-XmlOverRange $rng1 $body1

# Change 2: wrap WorkersRightsComputer.rc run pair with spellStart/spellEnd
$r2 = $d.Content
$r2.Find.Execute("WorkersRightsComputer.rc", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$rng2 = $d.Range($r2.Start, $r2.End)
$body2 = '<w:proofErr w:type="spellStart"/><w:r w:rsidRPr="005C7EB6"><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>WorkersRightsComputer</w:t></w:r><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>.rc</w:t></w:r><w:proofErr w:type="spellEnd"/>'
Insert-XmlOverRange $rng2 $body2

# Change 3: By binary search on files ... *.rc file. -> wrap "rc" with spellStart/spellEnd
$r3 = $d.Content
$r3.Find.Execute("By binary search on files – the problem was with the *.rc file.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$rng3 = $d.Range($r3.Start, $r3.End)
$body3 = '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>By binary search on files – the problem was with the *.</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>rc</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> file.</w:t></w:r>'
Insert-XmlOverRange $rng3 $body3

# Change 4: Set new version to "v1.3.0.beta1". -> wrap "0.beta" with gramStart/gramEnd
$r4 = $d.Content
$r4.Find.Execute("Set new version to “v1.3.0.beta1”.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$rng4 = $d.Range($r4.Start, $r4.End)
$body4 = '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>Set new version to “v1.3.</w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>0.beta</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>1”.</w:t></w:r>'
Insert-XmlOverRange $rng4 $body4

# Change 5: Change letter paragraph -> append "- " + highlighted "Done"
$r5 = $d.Content
$r5.Find.Execute("Change letter – so that not all letters will be with Iris’s details – high priority", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$rng5 = $d.Range($r5.Start, $r5.End)
$body5 = '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>Change letter – so that not all letters will be with Iris’s details – high priority</w:t></w:r><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">- </w:t></w:r><w:r><w:rPr><w:highlight w:val="yellow"/><w:lang w:val="en-US"/></w:rPr><w:t>Done</w:t></w:r>'
Insert-XmlOverRange $rng5 $body5

# Change 6: "...of contact letter" paragraph -> append "- " + highlighted "Done"
$r6 = $d.Content
$r6.Find.Execute("Allow users to select details of contact letter", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$rng6 = $d.Range($r6.Start, $r6.End)
$body6 = '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>Allow users to select details</w:t></w:r><w:r w:rsidR="006C4D47"><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> of contact letter</w:t></w:r><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">- </w:t></w:r><w:r><w:rPr><w:highlight w:val="yellow"/><w:lang w:val="en-US"/></w:rPr><w:t>Done</w:t></w:r>'
Insert-XmlOverRange $rng6 $body6

# Change 7: "On starting, check..." paragraph -> append " " + "- " + highlighted "Done"
$r7 = $d.Content
$r7.Find.Execute("On starting, check that save directory exists – if not force user to config.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$rng7 = $d.Range($r7.Start, $r7.End)
$body7 = '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>On starting, check that save directory exists – if not force user to config.</w:t></w:r><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">- </w:t></w:r><w:r><w:rPr><w:highlight w:val="yellow"/><w:lang w:val="en-US"/></w:rPr><w:t>Done</w:t></w:r>'
Insert-XmlOverRange $rng7 $body7

Write-Host "All changes applied."
